$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.444118618965149
$ws.Range("B1").Value = 2.730725526809692
$ws.Range("C1").Value = 3.263580083847046
$ws.Range("D1").Value = 3.126994371414185
$ws.Range("E1").Value = 2.373293399810791
